$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing date number formats onto the new date cells so we reuse the
# existing style entries (s="3" = m/d/yyyy, s="4" = d-mmm) instead of Excel
# inventing brand-new style/numFmt entries.
$ws.Cells.Item(11, 1).Copy()
$ws.Cells.Item(15, 1).PasteSpecial(-4122)
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$ws.Cells.Item(18, 1).PasteSpecial(-4122)

$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 15: 2/21/2018 - Generic Dao/Testing - 2 hours
$ws.Cells.Item(15, 1).Value = 43152
$ws.Cells.Item(15, 2).Value = "Generic Dao/Testing"
$ws.Cells.Item(15, 5).Value = 2

# Row 16: 2/22/2018 - Fixing tests/DB - 2 hours
$ws.Cells.Item(16, 1).Value = 43153
$ws.Cells.Item(16, 2).Value = "Fixing tests/DB"
$ws.Cells.Item(16, 5).Value = 2

# Row 17: 2/26/2018 - AWS Setup / First Deploy - 6 hours
$ws.Cells.Item(17, 1).Value = 43157
$ws.Cells.Item(17, 2).Value = "AWS Setup / First Deploy"
$ws.Cells.Item(17, 5).Value = 6

# Row 18: 2/27/2018 - no activity/hours logged yet
$ws.Cells.Item(18, 1).Value = 43158

# Leave selection where the author left off
$ws.Range("B18").Select()
